$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "Densidade"
$ws.Range("D2").Value = 2000
$ws.Range("E2").Value = 1000
$ws.Range("G2").Value = 0.99

# Remove rows 3 and 4 entirely (they are no longer part of the data)
$ws.Range("A3:H4").Delete()
